$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.228.87"
$ws.Range("E2").Value = "  -0.77%  "

$ws.Range("D3").Value = "1.926.51"
$ws.Range("E3").Value = "  -3.84%  "

$ws.Range("E4").Value = "  -0.09%  "

$ws.Range("D5").Value = "'241.70"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.04%  "

$ws.Range("D6").Value = "'0.604"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -4.52%  "

$ws.Range("E7").Value = "  -0.06%  "

$ws.Range("D8").Value = "'56.17"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -9.12%  "

$ws.Range("D9").Value = "'0.361"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -5.81%  "

$ws.Range("D10").Value = "'54.51"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -5.67%  "

$ws.Range("D11").Value = "'0.0822"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +5.55%  "

$ws.Range("E12").Value = "  -0.82%  "

$ws.Range("B13").Value = "Polygon"
$ws.Range("C13").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D13").Value = "'0.807"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -8.86%  "

$ws.Range("B14").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C14").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D14").Value = "2.198.83"
$ws.Range("E14").Value = "  -4.25%  "

$ws.Range("B15").Value = "Avalanche"
$ws.Range("C15").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D15").Value = "'20.78"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -9.42%  "

$ws.Range("B16").Value = "Chainlink"
$ws.Range("C16").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D16").Value = "'13.17"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -7.52%  "

$ws.Range("D17").Value = "'5.16"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -6.15%  "

$ws.Range("D18").Value = "1.912.37"
$ws.Range("E18").Value = "  -4.40%  "

$ws.Range("D19").Value = "36.100.93"
$ws.Range("E19").Value = "  -0.98%  "

$ws.Range("D20").Value = "'69.15"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.73%  "

$ws.Range("D21").Value = "0.0₃0853"
$ws.Range("E21").Value = "  -2.22%  "

$ws.Range("D22").Value = "'225.46"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.81%  "

$ws.Range("D23").Value = "'4.93"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -6.94%  "

$ws.Range("E24").Value = "  +0.06%  "

$ws.Range("D25").Value = "'2.42"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.08%  "

$ws.Range("D26").Value = "'2.25"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.51%  "

$ws.Range("D27").Value = "'9.24"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -6.21%  "

$ws.Range("D28").Value = "'162.52"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.88%  "

$ws.Range("D29").Value = "'19.09"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -5.19%  "

$ws.Range("D30").Value = "'0.117"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -14.82%  "

$ws.Range("E31").Value = "  -3.35%  "

$ws.Range("D32").Value = "'1.12"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.80%  "

$ws.Range("D33").Value = "'4.62"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -7.15%  "

$ws.Range("D34").Value = "'0.0616"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.47%  "

$ws.Range("D35").Value = "'4.22"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -5.04%  "

$ws.Range("E36").Value = "  -0.13%  "

$ws.Range("B37").Value = "THORChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D37").Value = "'5.95"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -8.80%  "

$ws.Range("B38").Value = "WEMIXToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D38").Value = "'1.78"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.68%  "

$ws.Range("D39").Value = "'2.13"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -9.60%  "

$ws.Range("D40").Value = "'2.85"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -11.23%  "

$ws.Range("D41").Value = "'0.0952"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -5.13%  "

$ws.Range("D42").Value = "'2.84"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.64%  "

$ws.Range("D43").Value = "'1.15"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -8.41%  "

$ws.Range("D44").Value = "'0.0206"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.60%  "

$ws.Range("B45").Value = "InjectiveProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D45").Value = "'15.43"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -7.17%  "

$ws.Range("B46").Value = "Maker"
$ws.Range("C46").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D46").Value = "1.333.00"
$ws.Range("E46").Value = "  -1.62%  "

$ws.Range("D47").Value = "'1.02"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -8.83%  "

$ws.Range("D48").Value = "'86.71"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -6.76%  "

$ws.Range("D49").Value = "'7.16"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -5.57%  "

$ws.Range("D50").Value = "'2.80"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.65%  "

$ws.Range("D51").Value = "'45.10"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.98%  "
